# Edit script for 09_Layout_Dashboard_Contable.xlsx
# - Update dashboard title on Dashboard_Layout sheet
# - Rewrite the step-by-step instructions on the Instrucciones sheet
#   (including several row-height tweaks)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Dashboard_Layout ---
$dashboard = $wb.Worksheets.Item("Dashboard_Layout")
$dashboard.Range("A1").Value = "PLANTILLA: Construye Tu Dashboard Aqui"

# --- Sheet 2: Instrucciones ---
$instrucciones = $wb.Worksheets.Item("Instrucciones")

$instrucciones.Range("B3").Value = "PLANTILLA DE TRABAJO -- Usa este archivo para construir tu dashboard durante la clase."

$instrucciones.Range("B4").Value = "PASO 1: Abre los archivos de ejercicios anteriores (Modulos 1-3) donde tienes tus datos de nomina."
$instrucciones.Rows.Item(4).RowHeight = 32

$instrucciones.Range("B5").Value = "PASO 2: Crea una Tabla Dinamica desde tus datos (Insertar > Tabla Dinamica) y pegala en esta hoja."
$instrucciones.Rows.Item(5).RowHeight = 32

$instrucciones.Range("B6").Value = "PASO 3: Inserta Segmentadores vinculados a tu TD (clic en TD > Insertar > Segmentacion de datos)."
$instrucciones.Rows.Item(6).RowHeight = 32

$instrucciones.Range("B7").Value = "PASO 4: Crea graficos desde tu TD y colocalos en las areas marcadas con bordes punteados."
$instrucciones.Rows.Item(7).RowHeight = 30

$instrucciones.Range("B8").Value = "PASO 5: Reemplaza '$0.00' en los KPIs con formulas =SUBTOTAL(109,...) que apunten a tu tabla."

$instrucciones.Range("B9").Value = "PASO 6: Oculta las lineas de cuadricula (Vista > desmarcar 'Lineas de cuadricula')."
$instrucciones.Rows.Item(9).RowHeight = 30

$instrucciones.Range("B10").Value = "Si quieres ver una solucion completa de referencia, abre el archivo: 10_Dashboard_Final_Integrado.xlsx"
$instrucciones.Rows.Item(10).RowHeight = 34

$instrucciones.Range("B11").Value = "Tip: Inmoviliza paneles en fila 5 para que los KPIs queden fijos al desplazarte."

$instrucciones.Range("B12").Value = "Este archivo es TU espacio de trabajo. El instructor trabaja en paralelo con el mismo template."
